$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "total_clp" header + values ----------------------------
$ws.Range("G1").Value = "total_clp"
$ws.Range("G2").Value = 20083803
$ws.Range("G3").Value = 20083803

# Style the new header cell like the other header cells (bold font, centered,
# top-aligned) but with a distinct border (thin left/right only, no top/bottom).
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G1").VerticalAlignment = -4160     # xlTop
$ws.Range("G1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$ws.Range("G1").Borders.Item(10).LineStyle = 1  # xlEdgeRight, xlContinuous

# --- Column F ("MB") values changed from 34 to 6 ---------------------------
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6

# --- Column B ("dedicacion") loses its explicit style ----------------------
$ws.Range("B2:B3").Style = "Normal"

# --- Selection moves to F3 --------------------------------------------------
$ws.Range("F3").Select() | Out-Null
